$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8628879189491272
$ws.Range("B1").Value = 1.319664716720581
$ws.Range("D1").Value = -1
$ws.Range("E1").Value = 1.555025339126587
